$d = $word.ActiveDocument

# Locate the paragraph that ends with the "LOQ4009: ..." requirement line.
# The three paragraphs that immediately follow it (a blank paragraph, the
# "Ver no Jupiter..." line, and the "(c) 2020 ..." footer line) are the ones
# being removed by this edit; the blank paragraph that precedes the
# page-break paragraph at the very end of the document must stay untouched.
$targetIndex = -1
$i = 1
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*LOQ4009*") {
        $targetIndex = $i
    }
    $i = $i + 1
}

if ($targetIndex -ge 1) {
    $startIndex = $targetIndex + 1
    $endIndex = $targetIndex + 3

    $rangeStart = $d.Paragraphs($startIndex).Range.Start
    $rangeEnd = $d.Paragraphs($endIndex).Range.End

    $r = $d.Range($rangeStart, $rangeEnd)
    $r.Delete()
}
